$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.823.88'
$ws.Range("E2").Value = '  -5.11%  '

# Row 3
$ws.Range("D3").Value = '2.491.29'
$ws.Range("E3").Value = '  -3.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.66%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.24%  '

# Row 9
$ws.Range("D9").Value = '2.519.28'
$ws.Range("E9").Value = '  -2.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.13%  '

# Row 11
$ws.Range("E11").Value = '  -2.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.25%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.34%  '

# Row 14
$ws.Range("D14").Value = '2.930.70'
$ws.Range("E14").Value = '  -3.40%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.46%  '

# Row 16
$ws.Range("D16").Value = '58.739.06'
$ws.Range("E16").Value = '  -5.02%  '

# Row 17
$ws.Range("E17").Value = '  -5.26%  '

# Row 18
$ws.Range("D18").Value = '2.517.04'
$ws.Range("E18").Value = '  -2.53%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.85%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

# Row 23
$ws.Range("E23").Value = '  -5.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.438'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -11.29%  '

# Row 26
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.621.20'
$ws.Range("E26").Value = '  -2.82%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.162'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.83%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.37%  '

# Row 29
$ws.Range("E29").Value = '  -5.61%  '

# Row 30
$ws.Range("E30").Value = '  -7.45%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0771'
$ws.Range("E31").Value = '  -8.49%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.41%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.28%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.66%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.01%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.10%  '

# Row 39
$ws.Range("E39").Value = '  -11.46%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.02%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '307.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.16%  '

# Row 43
$ws.Range("E43").Value = '  -7.25%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.803'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -12.70%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.34%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.595'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.95%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.37%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.13%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0926'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.42%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.29%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0515'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.18%  '
